$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the submitter email address and turn it into a mailto: hyperlink
# (AK2 held "sindhube@gmail.com" before; it becomes the new corporate
# address and gets styled as a hyperlink, same as a user finishing/fixing
# the "email" column in Excel).
$cell = $ws.Range("AK2")
$cell.Value = "sindhuja.e@mstsolutions.com"
$ws.Hyperlinks.Add($cell, "mailto:sindhuja.e@mstsolutions.com")

# Leave the grid scrolled over to the right with AD6 selected, matching
# where the user ended up after finishing the edit.
$ws.Range("AD6").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 29
$excel.ActiveWindow.ScrollRow = 1
